$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "Written by" -> "Written", "Reviewed by" -> "Reviewed"
$ws.Range("L1").Value = "Written"
$ws.Range("M1").Value = "Reviewed"

# Widen column B to (approximately) match column C's width (8.88671875, customWidth)
$ws.Columns.Item(2).ColumnWidth = 8

# Scroll the view so column K becomes the left-most visible column, and move
# the active selection to M16 (matches the saved view state in the workbook)
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
$ws.Range("M16").Select()
